$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Prova1")
$ws2 = $wb.Worksheets.Item("Prova2")

# --- "Nome" column (A3) for both sheets: "a" -> "Teste2" ---
$ws1.Range("A3").Value = "Teste2"
$ws2.Range("A3").Value = "Teste2"

# --- "Data de nascimento" column (B3) for both sheets: "a" -> "12/12/2012" ---
# Temporarily force a text number format so the date-like string is stored as a
# literal string instead of being auto-converted to a date serial value, then
# restore the cell's original look (border/centering, same as column A) by
# copying the formatting over from A3 (which keeps its original style) so no
# stray date-format style sticks to the cell.
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "12/12/2012"
$ws1.Range("A3").Copy()
$ws1.Range("B3").PasteSpecial(-4122)

$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "12/12/2012"
$ws2.Range("A3").Copy()
$ws2.Range("B3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Prova1 (sheet1) row 3 remaining updates ---
$ws1.Range("H3").Value = "00:00:11"
$ws1.Range("I3").Value = 0
$ws1.Range("J3").Value = 1

# --- Prova2 (sheet2) row 3 remaining updates ---
$ws2.Range("F3").Value = 100
$ws2.Range("G3").Value = 50
$ws2.Range("H3").Value = "00:00:08"
$ws2.Range("I3").Value = 0
$ws2.Range("J3").Value = 0.5
